# New weekly price report: insert a new record for Ciboulette (Femacal de
# La Calera) right before the existing row 248, shifting the remaining
# historical rows down by one (old row 349 becomes row 350).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 248..349 down to 249..350 and create a fresh (blank, but
# formatted like the row above) row 248.
$ws.Rows.Item(248).Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(248, 1).Value  = 3
$ws.Cells.Item(248, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(248, 3).Value  = "Coquimbo"
$ws.Cells.Item(248, 4).Value  = 44784
$ws.Cells.Item(248, 5).Value  = 5
$ws.Cells.Item(248, 6).Value  = 100112039
$ws.Cells.Item(248, 7).Value  = "Ciboulette"
$ws.Cells.Item(248, 8).Value  = "Sin especificar"
$ws.Cells.Item(248, 9).Value  = "Primera"
$ws.Cells.Item(248, 10).Value = 120
$ws.Cells.Item(248, 11).Value = 1500
$ws.Cells.Item(248, 12).Value = 1500
$ws.Cells.Item(248, 13).Value = 1500
$ws.Cells.Item(248, 14).Value = "`$/docena de atados"
$ws.Cells.Item(248, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(248, 16).Value = 500
$ws.Cells.Item(248, 17).Value = 3
$ws.Cells.Item(248, 18).Value = "Hortaliza"
